$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "Table"

$ws.Range("A1").Value = "Numbers"
$ws.Range("A2").Value = 5
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 1

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:A7"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.ShowTotals = $true
$lo.ListColumns.Item(1).TotalsCalculation = 1
Write-Output "done"
